$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert two new paragraphs after "...you through the process of resolving
#    that merge conflict." -- one blank paragraph, then one paragraph with
#    "Remember to work within the KitClient when working locally." (with
#    proofErr spell-check markers bracketing "KitClient").
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("you through the process of resolving that merge conflict.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find merge-conflict intro sentence" }

$point = $d.Range($rng.End, $rng.End)
$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Remember to work within the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>KitClient</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> when working locally.</w:t></w:r></w:p>
'@
$point.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: remove it from in front of
#    "ii. Highlight lines that cause conflicting changes in red." and add it
#    in front of "a. " at the start of the "Use the highlight tool..."
#    paragraph instead.
# ---------------------------------------------------------------------------

# 2a. Remove it from the "ii. Highlight..." paragraph.
$rngIi = $d.Content
$foundIi = $rngIi.Find.Execute("Highlight lines that cause conflicting changes in red.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundIi) { throw "Could not find 'ii. Highlight...red.' paragraph text" }
$paraIiFull = $d.Range($rngIi.Start, $rngIi.End + 1)
$paraIiXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:tab/><w:t xml:space="preserve">ii. Highlight lines that </w:t></w:r><w:r w:rsidR="008000F4"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">cause </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>conflicting changes in red.</w:t></w:r></w:p>
'@
$paraIiFull.InsertXML($paraIiXml)

# 2b. Add it to the "a. Use the highlight tool..." paragraph.
$rngA = $d.Content
$foundA = $rngA.Find.Execute("a. Use the highlight tool to mark the lines in the feature branch and/or the main branch above as indicated below. Use the example in the slides as a guide for the highlighting.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundA) { throw "Could not find 'a. Use the highlight tool...' paragraph text" }
$paraAFull = $d.Range($rngA.Start, $rngA.End + 1)
$paraAXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">a. </w:t></w:r><w:r w:rsidR="001D6C29"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Use the highlight tool to mark the lines in the feature branch and/or the </w:t></w:r><w:r w:rsidR="001D6C29" w:rsidRPr="008000F4"><w:rPr><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:cstheme="minorHAnsi"/></w:rPr><w:t>main</w:t></w:r><w:r w:rsidR="001D6C29"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> branch above as </w:t></w:r><w:r w:rsidR="00D42A48"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">indicated below. </w:t></w:r><w:r w:rsidR="009A1DE6"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>Use the example in the slides as a guide for the highlighting.</w:t></w:r></w:p>
'@
$paraAFull.InsertXML($paraAXml)

# ---------------------------------------------------------------------------
# 3) Remove the gramStart/gramEnd proofErr markers around "customized" and
#    merge the three runs they split into a single run.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("20. Complete the table below", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find '20. Complete the table below...'" }

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Note that the tasks listed are in approximately the same order as they appear in this activity.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find 'Note that the tasks listed...' sentence" }

$paraCustomizedFull = $d.Range($rng1.Start, $rng2.End + 1)
$paraCustomizedXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>20</w:t></w:r><w:r w:rsidR="00D47D92"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">. Complete the table below by filling in the right-hand column with the commands that accomplish the task listed on the left.  Use the </w:t></w:r><w:r w:rsidR="00D47D92" w:rsidRPr="00A34EDE"><w:rPr><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:cstheme="minorHAnsi"/></w:rPr><w:t>&lt;&#8230;&gt;</w:t></w:r><w:r w:rsidR="00D47D92"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> notation appropriately to indicate parameters that need to customized for each use.  Note that the tasks listed are in approximately the same order as they appear in this activity.</w:t></w:r></w:p>
'@
$paraCustomizedFull.InsertXML($paraCustomizedXml)

Write-Output "All edits applied"
